$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.01586292554596041
$ws.Range("C2").Value = 0.6854252723133342
$ws.Range("D2").Value = 0.4690748481128704
$ws.Range("E2").Value = 0.5942751998756042
$ws.Range("F2").Value = 0.09260064950199673
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = 0.1762291725437023
$ws.Range("J2").Value = 0.5079728618234201
$ws.Range("B3").Value = 0.4469482919323776
$ws.Range("C3").Value = 0.4293442456321315
$ws.Range("D3").Value = 0.5597239046618772
$ws.Range("E3").Value = 0.3941813052359234
$ws.Range("F3").Value = 0.5155554027454251
$ws.Range("I3").Value = 0.3754703218770988
$ws.Range("J3").Value = 0.4899234134726962
$ws.Range("B4").Value = 0.1498073400962025
$ws.Range("C4").Value = 0.2456447693183586
$ws.Range("D4").Value = 0.179617747755455
$ws.Range("E4").Value = 0.1650435633692065
$ws.Range("F4").Value = 0.2192535997906139
$ws.Range("G4").Value = -0.04046440566198269
$ws.Range("H4").Value = -0.05227215521641281
$ws.Range("I4").Value = 0.5352968665696594
$ws.Range("J4").Value = 0.1743013859539189
$ws.Range("B5").Value = 0.9230451578990424
$ws.Range("C5").Value = 0.8490787306869069
$ws.Range("D5").Value = 0.9230451578990424
$ws.Range("E5").Value = 0.7780876219875478
$ws.Range("F5").Value = 0.6779570593759783
$ws.Range("I5").Value = 0.9230451578990424
$ws.Range("J5").Value = 0.3674804716401385
$ws.Range("B6").Value = 0.3872850122850123
$ws.Range("C6").Value = 0.3137373432068913
$ws.Range("D6").Value = 0.4954701655732584
$ws.Range("E6").Value = 0.2873527588344699
$ws.Range("F6").Value = 0.06943244021895703
$ws.Range("G6").Value = 0.05688789737869489
$ws.Range("H6").Value = 0.01624896722665934
$ws.Range("I6").Value = 0.2380049584366341
$ws.Range("J6").Value = 0.3884246363354998
$ws.Range("B7").Value = 0.3581887149080505
$ws.Range("C7").Value = 0.6834901485635159
$ws.Range("D7").Value = 0.3581887149080505
$ws.Range("E7").Value = 0.5093923760575323
$ws.Range("F7").Value = 0.0358622987011715
$ws.Range("G7").Value = -0.04017794754739839
$ws.Range("I7").Value = 0.05696576151121614
$ws.Range("J7").Value = -0.05598121308442266
$ws.Range("B8").Value = 0.774672769640816
$ws.Range("C8").Value = 0.6401957566748393
$ws.Range("D8").Value = 0.292630114277278
$ws.Range("E8").Value = 0.4461260139534836
$ws.Range("F8").Value = 0.4993337233216952
$ws.Range("G8").Value = 0.015030426365802
$ws.Range("I8").Value = 0.4910031696332042
$ws.Range("J8").Value = 0.06305699008272014
$ws.Range("B9").Value = -0.07339844711044838
$ws.Range("C9").Value = 0.8701742752346809
$ws.Range("D9").Value = 0.3074685757552467
$ws.Range("E9").Value = 0.5009917975660751
$ws.Range("F9").Value = 0.3131865223164396
$ws.Range("G9").Value = 0.1178150233837516
$ws.Range("H9").Value = 0.2161308124973622
$ws.Range("I9").Value = 0.3724218022367568
$ws.Range("J9").Value = 0.05988086261185167
$ws.Range("B10").Value = 0.02222222222222222
$ws.Range("C10").Value = -0.02401182120428521
$ws.Range("D10").Value = -0.03455708025001789
$ws.Range("E10").Value = 0.2259847961299239
$ws.Range("F10").Value = -0.01587301587301587
$ws.Range("G10").Value = -0.1065903603684166
$ws.Range("I10").Value = 0.1200312877764346
$ws.Range("J10").Value = -0.0772072404876247
$ws.Range("B11").Value = 0.09929485733413863
$ws.Range("C11").Value = 0.0291724293203821
$ws.Range("D11").Value = 0.06624007900101087
$ws.Range("E11").Value = 0.01536751758084563
$ws.Range("F11").Value = 0.06416637225208294
$ws.Range("G11").Value = 0.01917226442420983
$ws.Range("H11").Value = 0.02864970942634217
$ws.Range("I11").Value = 0.1465098535459739
$ws.Range("J11").Value = 0.02089691122409309
$ws.Range("B12").Value = 0.01451679800912478
$ws.Range("D12").Value = -0.01503759398496235
$ws.Range("E12").Value = 0.01451679800912478
$ws.Range("F12").Value = -0.01027960526315787
$ws.Range("I12").Value = -0.02932551319648097
$ws.Range("J12").Value = -0.02704987320371931
$ws.Range("B13").Value = -0.06478873239436632
$ws.Range("C13").Value = 0.07046070460704609
$ws.Range("D13").Value = -0.06944444444444441
$ws.Range("E13").Value = -0.02997275204359669
$ws.Range("F13").Value = -0.05538461538461534
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = -0.0576923076923077
$ws.Range("J13").Value = 0.1812865497076024
$ws.Range("B14").Value = 0.3769230769230769
$ws.Range("C14").Value = -0.09197080291970815
$ws.Range("D14").Value = -0.02307692307692308
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0.1440922190201729
$ws.Range("G14").Value = 0.1991999999999999
$ws.Range("I14").Value = 0.04106628242074932
$ws.Range("J14").Value = 0.04676753782668504
$ws.Range("B15").Value = 0.5130472475075856
$ws.Range("C15").Value = 0.3438428366075707
$ws.Range("D15").Value = 0.4487714418173389
$ws.Range("E15").Value = 0.3376476145488899
$ws.Range("F15").Value = 0.731333640978311
$ws.Range("G15").Value = 0.1270668176670441
$ws.Range("H15").Value = -0.0006139677666922268
$ws.Range("I15").Value = 0.6088208820882087
$ws.Range("J15").Value = 0.3533623609095307
$ws.Range("B16").Value = 0.41771662210604
$ws.Range("C16").Value = 0.3733871127453246
$ws.Range("D16").Value = 0.4856464407187964
$ws.Range("E16").Value = 0.4238281244230387
$ws.Range("F16").Value = 0.3433221522680295
$ws.Range("G16").Value = 0.07218113831635374
$ws.Range("H16").Value = 0.05756422450970444
$ws.Range("I16").Value = 0.4449905980697924
$ws.Range("J16").Value = 0.291070299888343
$ws.Range("B17").Value = 0.2100437297221046
$ws.Range("C17").Value = -0.009026994184532606
$ws.Range("D17").Value = 0.1449934068400528
$ws.Range("E17").Value = 0.07204780742361756
$ws.Range("F17").Value = 0.2100437297221046
$ws.Range("G17").Value = 0.1032793334048564
$ws.Range("I17").Value = 0.2153565924057728
$ws.Range("J17").Value = -0.09471214935001485
$ws.Range("B18").Value = 0.0229667881006179
$ws.Range("C18").Value = 0.01858277327451917
$ws.Range("D18").Value = 0.03845373496969431
$ws.Range("E18").Value = 0.0859629674552513
$ws.Range("F18").Value = 0.09809833824115108
$ws.Range("G18").Value = -0.007520582002594521
$ws.Range("H18").Value = 0.003400311055631334
$ws.Range("I18").Value = 0.02891295783279457
$ws.Range("J18").Value = 0.05517714813032964
$ws.Range("B19").Value = 0.4277524218092092
$ws.Range("C19").Value = 0.3865853726937875
$ws.Range("D19").Value = 0.4038478214192654
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0.4855928333401088
$ws.Range("G19").Value = 0.04811294912513467
$ws.Range("H19").Value = 0.1141997913488729
$ws.Range("I19").Value = 0.3497656062829619
$ws.Range("J19").Value = 0.2381696085989845
$ws.Range("B20").Value = 0.3259359903381642
$ws.Range("C20").Value = 0.2378378378378379
$ws.Range("D20").Value = 0.1561557030778516
$ws.Range("F20").Value = 0.01201201201201204
$ws.Range("I20").Value = 0.03083716926761719
$ws.Range("J20").Value = 0.1259345531051584
$ws.Range("B21").Value = -0.0003192649244378704
$ws.Range("C21").Value = -0.000452379759883236
$ws.Range("D21").Value = 0.1276766954267066
$ws.Range("E21").Value = 0.02582549127837833
$ws.Range("F21").Value = 0.03680481558334739
$ws.Range("I21").Value = 0.5409740077900257
$ws.Range("J21").Value = -0.1070392828895331
$ws.Range("B22").Value = 0.1805643922215898
$ws.Range("C22").Value = 0.1547785411539182
$ws.Range("D22").Value = 0.3644577333044043
$ws.Range("E22").Value = -0.006366041680322763
$ws.Range("F22").Value = 0.3044498338916452
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0.03069544202482084
$ws.Range("J22").Value = 0.1145051826490301
